$d = $word.ActiveDocument

# 1. "<closing quote>." -> "<closing quote>" in the "-Optimización ... admitirpalabra"." paragraph
$d.Content.Find.Execute("”.", $true, $false, $false, $false, $false, $true, 1, $false, "”", 2)

# 2. Replace the "-No se pide..." paragraph text with the new sentence.
$d.Content.Find.Execute("-No se pide la palabra a la hora de jugar. Solo la posición y la orientación.", $true, $false, $false, $false, $false, $true, 1, $false, "-Palabra contiene espacio: se vuelve a pedir.", 2)

# 3. Move the "_GoBack" bookmark so it sits right after "Diagrama:" (collapsed, at paragraph end)
#    instead of at its original location at the end of the document. Bookmarks.Add placed at the
#    exact end-of-paragraph offset misbehaves, so a temporary marker character is used to get a
#    safe, reliable offset, then removed again.
$diagramaRange = $d.Content
$diagramaRange.Find.Execute("Diagrama:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$diagramaRange.Collapse(0)
$bookmarkPos = $diagramaRange.Start
$diagramaRange.InsertAfter("@@MARK@@")

$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$markerRange = $d.Content
$markerRange.Find.Execute("@@MARK@@", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerRange.Delete()

# 4. Remove the "Dibujar..." / "Jugar..." / empty paragraphs that used to sit between
#    "Diagrama:" and "Caja Negra:".
$startRange = $d.Content
$startRange.Find.Execute("Dibujar", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endRange = $d.Content
$endRange.Find.Execute("Caja Negra:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$toDelete = $d.Range($startRange.Start, $endRange.Start)
$toDelete.Delete()
